$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the three new variables
$ws.Range("BG1").Value = "dist_trav_20min_body_out"
$ws.Range("BH1").Value = "dist_trav_25min_body_out"
$ws.Range("BI1").Value = "dist_trav_30min_body_out"

# Populate data rows 2-78 for BG, BH, BI columns
$ws.Range("BG2").Value = 0
$ws.Range("BH2").Value = 0
$ws.Range("BI2").Value = 0
$ws.Range("BG3").Value = 0
$ws.Range("BH3").Value = 0
$ws.Range("BI3").Value = 0
$ws.Range("BG4").Value = 29.065045753
$ws.Range("BH4").Value = 43.707970517
$ws.Range("BI4").Value = 53.959036928
$ws.Range("BG5").Value = 51.8091856353
$ws.Range("BH5").Value = 0
$ws.Range("BI5").Value = 0
$ws.Range("BG6").Value = 0
$ws.Range("BH6").Value = 0
$ws.Range("BI6").Value = 0
$ws.Range("BG7").Value = 0
$ws.Range("BH7").Value = 0
$ws.Range("BI7").Value = 0
$ws.Range("BG8").Value = 0
$ws.Range("BH8").Value = 0
$ws.Range("BI8").Value = 0
$ws.Range("BG9").Value = 0
$ws.Range("BH9").Value = 0
$ws.Range("BI9").Value = 0
$ws.Range("BG10").Value = 97.496614224
$ws.Range("BH10").Value = 116.082207546
$ws.Range("BI10").Value = 134.850822292
$ws.Range("BG11").Value = 0
$ws.Range("BH11").Value = 0
$ws.Range("BI11").Value = 0
$ws.Range("BG12").Value = 0
$ws.Range("BH12").Value = 0
$ws.Range("BI12").Value = 0
$ws.Range("BG13").Value = 55.8808590596
$ws.Range("BH13").Value = 69.97903018229999
$ws.Range("BI13").Value = 87.1682587289
$ws.Range("BG14").Value = 0
$ws.Range("BH14").Value = 0
$ws.Range("BI14").Value = 0
$ws.Range("BG15").Value = 0
$ws.Range("BH15").Value = 0
$ws.Range("BI15").Value = 0
$ws.Range("BG16").Value = 0
$ws.Range("BH16").Value = 0
$ws.Range("BI16").Value = 0
$ws.Range("BG17").Value = 71.51811290969999
$ws.Range("BH17").Value = 83.4898792445
$ws.Range("BI17").Value = 94.83043631629999
$ws.Range("BG18").Value = 74.51979790576
$ws.Range("BH18").Value = 86.63974060347
$ws.Range("BI18").Value = 0
$ws.Range("BG19").Value = 0
$ws.Range("BH19").Value = 0
$ws.Range("BI19").Value = 0
$ws.Range("BG20").Value = 0
$ws.Range("BH20").Value = 0
$ws.Range("BI20").Value = 0
$ws.Range("BG21").Value = 0
$ws.Range("BH21").Value = 0
$ws.Range("BI21").Value = 0
$ws.Range("BG22").Value = 60.851160999
$ws.Range("BH22").Value = 75.42470776499999
$ws.Range("BI22").Value = 85.755232096
$ws.Range("BG23").Value = 113.053870372
$ws.Range("BH23").Value = 136.9123724667
$ws.Range("BI23").Value = 156.0485786097
$ws.Range("BG24").Value = 75.196236946
$ws.Range("BH24").Value = 89.062112793
$ws.Range("BI24").Value = 106.651700849
$ws.Range("BG25").Value = 0
$ws.Range("BH25").Value = 0
$ws.Range("BI25").Value = 0
$ws.Range("BG26").Value = 0
$ws.Range("BH26").Value = 0
$ws.Range("BI26").Value = 0
$ws.Range("BG27").Value = 0
$ws.Range("BH27").Value = 0
$ws.Range("BI27").Value = 0
$ws.Range("BG28").Value = 75.513032749
$ws.Range("BH28").Value = 94.486585324
$ws.Range("BI28").Value = 115.689266338
$ws.Range("BG29").Value = 43.7300336234
$ws.Range("BH29").Value = 68.1845941214
$ws.Range("BI29").Value = 85.9006444444
$ws.Range("BG30").Value = 57.8923168888
$ws.Range("BH30").Value = 79.8039617264
$ws.Range("BI30").Value = 97.59974628339999
$ws.Range("BG31").Value = 80.48688252620001
$ws.Range("BH31").Value = 81.5736840917
$ws.Range("BI31").Value = 83.6186979527
$ws.Range("BG32").Value = 81.55335022200001
$ws.Range("BH32").Value = 93.657725305
$ws.Range("BI32").Value = 109.372537696
$ws.Range("BG33").Value = 0
$ws.Range("BH33").Value = 0
$ws.Range("BI33").Value = 0
$ws.Range("BG34").Value = 106.71652432
$ws.Range("BH34").Value = 129.559830121
$ws.Range("BI34").Value = 155.295083402
$ws.Range("BG35").Value = 0
$ws.Range("BH35").Value = 0
$ws.Range("BI35").Value = 0
$ws.Range("BG36").Value = 117.346722634
$ws.Range("BH36").Value = 134.538752376
$ws.Range("BI36").Value = 0
$ws.Range("BG37").Value = 0
$ws.Range("BH37").Value = 0
$ws.Range("BI37").Value = 0
$ws.Range("BG38").Value = 82.65475698100001
$ws.Range("BH38").Value = 107.271065678
$ws.Range("BI38").Value = 130.592298752
$ws.Range("BG39").Value = 91.5006471005
$ws.Range("BH39").Value = 109.0869068761
$ws.Range("BI39").Value = 127.5684101063
$ws.Range("BG40").Value = 55.375485867
$ws.Range("BH40").Value = 55.375485867
$ws.Range("BI40").Value = 55.375485867
$ws.Range("BG41").Value = 0
$ws.Range("BH41").Value = 0
$ws.Range("BI41").Value = 0
$ws.Range("BG42").Value = 79.102341307
$ws.Range("BH42").Value = 105.3188281953
$ws.Range("BI42").Value = 0
$ws.Range("BG43").Value = 0
$ws.Range("BH43").Value = 0
$ws.Range("BI43").Value = 0
$ws.Range("BG44").Value = 71.281038159
$ws.Range("BH44").Value = 88.412876766
$ws.Range("BI44").Value = 103.163733472
$ws.Range("BG45").Value = 52.705811749
$ws.Range("BH45").Value = 63.016966618
$ws.Range("BI45").Value = 73.402168129
$ws.Range("BG46").Value = 51.872209149
$ws.Range("BH46").Value = 66.26112639500001
$ws.Range("BI46").Value = 80.42873317900001
$ws.Range("BG47").Value = 75.556626194
$ws.Range("BH47").Value = 94.40515643099999
$ws.Range("BI47").Value = 109.044053422
$ws.Range("BG48").Value = 97.624565797
$ws.Range("BH48").Value = 113.303162893
$ws.Range("BI48").Value = 137.560560446
$ws.Range("BG49").Value = 0
$ws.Range("BH49").Value = 0
$ws.Range("BI49").Value = 0
$ws.Range("BG50").Value = 72.61480394500001
$ws.Range("BH50").Value = 0
$ws.Range("BI50").Value = 0
$ws.Range("BG51").Value = 60.252062814
$ws.Range("BH51").Value = 76.023184211
$ws.Range("BI51").Value = 0
$ws.Range("BG52").Value = 83.365981737
$ws.Range("BH52").Value = 101.725452821
$ws.Range("BI52").Value = 0
$ws.Range("BG53").Value = 71.552019311
$ws.Range("BH53").Value = 92.298918266
$ws.Range("BI53").Value = 0
$ws.Range("BG54").Value = 93.87450332909999
$ws.Range("BH54").Value = 115.232602715
$ws.Range("BI54").Value = 128.1463967038
$ws.Range("BG55").Value = 75.85890205699999
$ws.Range("BH55").Value = 0
$ws.Range("BI55").Value = 0
$ws.Range("BG56").Value = 118.687446339
$ws.Range("BH56").Value = 135.011023018
$ws.Range("BI56").Value = 151.616284415
$ws.Range("BG57").Value = 0
$ws.Range("BH57").Value = 0
$ws.Range("BI57").Value = 0
$ws.Range("BG58").Value = 52.447994867
$ws.Range("BH58").Value = 0
$ws.Range("BI58").Value = 0
$ws.Range("BG59").Value = 85.908890035
$ws.Range("BH59").Value = 104.010235133
$ws.Range("BI59").Value = 0
$ws.Range("BG60").Value = 0
$ws.Range("BH60").Value = 0
$ws.Range("BI60").Value = 0
$ws.Range("BG61").Value = 0
$ws.Range("BH61").Value = 0
$ws.Range("BI61").Value = 0
$ws.Range("BG62").Value = 85.88790569299999
$ws.Range("BH62").Value = 95.810562464
$ws.Range("BI62").Value = 99.15457831800001
$ws.Range("BG63").Value = 74.8552253232
$ws.Range("BH63").Value = 91.9242211235
$ws.Range("BI63").Value = 107.9400711688
$ws.Range("BG64").Value = 0
$ws.Range("BH64").Value = 0
$ws.Range("BI64").Value = 0
$ws.Range("BG65").Value = 58.623882499
$ws.Range("BH65").Value = 76.571656111
$ws.Range("BI65").Value = 93.370780098
$ws.Range("BG66").Value = 0
$ws.Range("BH66").Value = 0
$ws.Range("BI66").Value = 0
$ws.Range("BG67").Value = 64.94725393
$ws.Range("BH67").Value = 80.330597814
$ws.Range("BI67").Value = 94.596163937
$ws.Range("BG68").Value = 0
$ws.Range("BH68").Value = 0
$ws.Range("BI68").Value = 0
$ws.Range("BG69").Value = 127.480112912
$ws.Range("BH69").Value = 162.091526286
$ws.Range("BI69").Value = 194.844537196
$ws.Range("BG70").Value = 70.9797245953
$ws.Range("BH70").Value = 0
$ws.Range("BI70").Value = 0
$ws.Range("BG71").Value = 95.82696878260001
$ws.Range("BH71").Value = 119.295739353
$ws.Range("BI71").Value = 141.0301472726
$ws.Range("BG72").Value = 0
$ws.Range("BH72").Value = 0
$ws.Range("BI72").Value = 0
$ws.Range("BG73").Value = 0
$ws.Range("BH73").Value = 0
$ws.Range("BI73").Value = 0
$ws.Range("BG74").Value = 101.916427861
$ws.Range("BH74").Value = 107.917648646
$ws.Range("BI74").Value = 127.79899625
$ws.Range("BG75").Value = 101.903279745
$ws.Range("BH75").Value = 123.107824233
$ws.Range("BI75").Value = 137.905860639
$ws.Range("BG76").Value = 0
$ws.Range("BH76").Value = 0
$ws.Range("BI76").Value = 0
$ws.Range("BG77").Value = 0
$ws.Range("BH77").Value = 0
$ws.Range("BI77").Value = 0
$ws.Range("BG78").Value = 106.924648905
$ws.Range("BH78").Value = 138.208306918
$ws.Range("BI78").Value = 158.930296279
